# Fruta / hortaliza, semanal
# Update Fecha (D), Calidad (I), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for rows 2..20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  Date="2021-07-14"; Calidad="Primera";  Volumen=400; PMin=15000; PMax=15000; PProm=15000; PKg=833 },
    @{ Row=3;  Date="2021-02-08"; Calidad="Primera";  Volumen=400; PMin=13000; PMax=13000; PProm=13000; PKg=722 },
    @{ Row=4;  Date="2021-02-08"; Calidad="Segunda";  Volumen=200; PMin=11000; PMax=11000; PProm=11000; PKg=611 },
    @{ Row=5;  Date="2021-02-08"; Calidad="Tercera";  Volumen=100; PMin=9000;  PMax=9000;  PProm=9000;  PKg=500 },
    @{ Row=6;  Date="2021-02-11"; Calidad="Primera";  Volumen=300; PMin=12000; PMax=12000; PProm=12000; PKg=667 },
    @{ Row=7;  Date="2021-02-11"; Calidad="Segunda";  Volumen=200; PMin=10000; PMax=10000; PProm=10000; PKg=556 },
    @{ Row=8;  Date="2021-02-11"; Calidad="Tercera";  Volumen=50;  PMin=8000;  PMax=8000;  PProm=8000;  PKg=444 },
    @{ Row=9;  Date="2021-02-22"; Calidad="Primera";  Volumen=400; PMin=12000; PMax=12000; PProm=12000; PKg=667 },
    @{ Row=10; Date="2021-02-22"; Calidad="Segunda";  Volumen=200; PMin=10000; PMax=10000; PProm=10000; PKg=556 },
    @{ Row=11; Date="2021-02-18"; Calidad="Primera";  Volumen=300; PMin=12000; PMax=12000; PProm=12000; PKg=667 },
    @{ Row=12; Date="2021-02-18"; Calidad="Segunda";  Volumen=200; PMin=10000; PMax=10000; PProm=10000; PKg=556 },
    @{ Row=13; Date="2021-02-16"; Calidad="Especial"; Volumen=300; PMin=12000; PMax=12000; PProm=12000; PKg=667 },
    @{ Row=14; Date="2021-02-16"; Calidad="Primera";  Volumen=300; PMin=10000; PMax=10000; PProm=10000; PKg=556 },
    @{ Row=15; Date="2021-02-16"; Calidad="Segunda";  Volumen=150; PMin=8000;  PMax=8000;  PProm=8000;  PKg=444 },
    @{ Row=16; Date="2021-07-19"; Calidad="Primera";  Volumen=250; PMin=15000; PMax=15000; PProm=15000; PKg=833 },
    @{ Row=17; Date="2021-07-19"; Calidad="Segunda";  Volumen=150; PMin=12000; PMax=12000; PProm=12000; PKg=667 },
    @{ Row=18; Date="2021-07-06"; Calidad="Primera";  Volumen=300; PMin=16000; PMax=16000; PProm=16000; PKg=889 },
    @{ Row=19; Date="2021-07-06"; Calidad="Segunda";  Volumen=200; PMin=12000; PMax=12000; PProm=12000; PKg=667 },
    @{ Row=20; Date="2021-02-02"; Calidad="Primera";  Volumen=200; PMin=15000; PMax=15000; PProm=15000; PKg=833 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = [DateTime]$r.Date
    $ws.Cells.Item($r.Row, 9).Value  = $r.Calidad
    $ws.Cells.Item($r.Row, 10).Value = $r.Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.PMin
    $ws.Cells.Item($r.Row, 12).Value = $r.PMax
    $ws.Cells.Item($r.Row, 13).Value = $r.PProm
    $ws.Cells.Item($r.Row, 16).Value = $r.PKg
}
